# aggiornamento fino a 9 agosto 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(329, 44403, 6, 29, 87.92408210290149),
    @(330, 44404, 1, 29, 87.92408210290149),
    @(331, 44405, 1, 30, 90.95594700300154),
    @(332, 44406, 1, 26, 78.82848740260134),
    @(333, 44407, 4, 22, 66.70102780220114),
    @(334, 44408, 3, 21, 63.66916290210109),
    @(335, 44409, 4, 20, 60.63729800200103),
    @(336, 44410, 1, 15, 45.47797350150077),
    @(337, 44411, 0, 14, 42.44610860140072),
    @(338, 44412, 0, 13, 39.41424370130067),
    @(339, 44413, 10, 22, 66.70102780220114),
    @(340, 44414, 0, 18, 54.57356820180092),
    @(341, 44415, 5, 20, 60.63729800200103),
    @(342, 44416, 0, 16, 48.50983840160082),
    @(343, 44417, 3, 18, 54.57356820180092)
)

# The date column (A) uses a special cell style (center/top aligned,
# bordered, custom date number format). Copy that formatting down from the
# last existing data row (A328) to each newly-added date cell so the new
# rows keep the same style index instead of minting a new one.
$ws.Range("A328").Copy() | Out-Null

foreach ($rowData in $data) {
    $r = $rowData[0]
    $ws.Cells.Item($r, 1).Value = $rowData[1]
    $ws.Cells.Item($r, 2).Value = $rowData[2]
    $ws.Cells.Item($r, 3).Value = $rowData[3]
    $ws.Cells.Item($r, 4).Value = $rowData[4]
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
